# FAST_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer footnote
#    from 2021-04-22 to 2021-04-23
#  - refresh Weight (D) and Percent Change (E) figures for rows 2-10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect so the cells below can be edited, then
# restore protection afterwards.
$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

# --- Disclaimer footnote text (cell A13) -----------------------------------
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) figures, rows 2-10 --------------------
$ws.Range("D2").Value = 0.1023795846995838
$ws.Range("E2").Value = 0.01260385463177527

$ws.Range("D3").Value = 0.1088371766601758
$ws.Range("E3").Value = 0.0160037380994098

$ws.Range("D4").Value = 0.1174928845136071
$ws.Range("E4").Value = 0.008698402656966753

$ws.Range("D5").Value = 0.1377284928063982
$ws.Range("E5").Value = 0.01215193948173199

$ws.Range("D6").Value = 0.1351501829778466
$ws.Range("E6").Value = 0.004576976421636791

$ws.Range("D7").Value = 0.1417793308301158
$ws.Range("E7").Value = 0.00704973442781287

$ws.Range("D8").Value = 0.1276897429188753
$ws.Range("E8").Value = 0.02466236054022319

$ws.Range("D9").Value = 0.1289426045933973
$ws.Range("E9").Value = 0.009840691742742935

$ws.Range("E10").Value = 0.0117639485339236

# Restore sheet protection to match the original state.
if ($wasProtected) {
    $ws.Protect()
}
